$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# D5: "Create a scrum master excel report" -> "Non development Task #83"
$ws.Range("D5").Value = "Non development Task #83"

# D8 was empty; give it the same formatting as D5/D15 (yellow task box) and a value
$ws.Range("D5").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = "Non development Task #85"

# D15: "Create a product owner excel report" -> "Non development Task #84"
$ws.Range("D15").Value = "Non development Task #84"

# New row 19, D19 formatted like A2 (orange header box) with its own value
$ws.Range("A2").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value = "Non development Task #86"

# Final selection left on D8, as in the saved workbook
$ws.Range("D8").Select()
